$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 10676.906
$ws.Range("I40").Value = 2109.75
$ws.Range("J40").Value = 13532.625
$ws.Range("K40").Value = 2109.75
$ws.Range("L40").Value = 13532.625
$ws.Range("M40").Value = -1934.75
$ws.Range("N40").Value = -13882.625

$ws.Range("H41").Value = 772.8125
$ws.Range("J41").Value = 730.6
$ws.Range("L41").Value = 730.6
$ws.Range("N41").Value = -1610.6

$ws.Range("H70").Value = 72982.07000000001
$ws.Range("I70").Value = 1050
$ws.Range("K70").Value = 3150
$ws.Range("M70").Value = -2880

$ws.Range("H73").Value = 72982.07000000001
$ws.Range("I73").Value = 1050
$ws.Range("K73").Value = 3150
$ws.Range("M73").Value = -2214

$ws.Range("H88").Value = 2751
$ws.Range("I88").Value = 2002.5
$ws.Range("J88").Value = 3499.5
$ws.Range("K88").Value = 2002.5
$ws.Range("L88").Value = 3499.5
$ws.Range("M88").Value = -1596.5
$ws.Range("N88").Value = -4311.5

$ws.Range("H91").Value = 2751
$ws.Range("I91").Value = 2002.5
$ws.Range("J91").Value = 3499.5
$ws.Range("K91").Value = 2002.5
$ws.Range("L91").Value = 3499.5
$ws.Range("M91").Value = -598.5
$ws.Range("N91").Value = -6307.5

$ws.Range("H93").Value = 49999.5
$ws.Range("J93").Value = 49999.5
$ws.Range("L93").Value = 49999.5
$ws.Range("N93").Value = -54991.5

$ws.Range("H101").Value = 525.3333
$ws.Range("I101").Value = 525.3333
$ws.Range("K101").Value = 1575.9999
$ws.Range("M101").Value = 46.00009999999997

$ws.Range("H103").Value = 1086.4286
$ws.Range("I103").Value = 1119.2
$ws.Range("J103").Value = 1004.5
$ws.Range("K103").Value = 3357.6
$ws.Range("L103").Value = 3013.5
$ws.Range("M103").Value = -2771.6
$ws.Range("N103").Value = -4185.5

$ws.Range("H113").Value = 4019.8
$ws.Range("I113").Value = 4019.8
$ws.Range("K113").Value = 4019.8
$ws.Range("M113").Value = -765.8000000000002

$ws.Range("H125").Value = 3211.1
$ws.Range("I125").Value = 3230.4285
$ws.Range("K125").Value = 29073.8565
$ws.Range("M125").Value = -26613.8565

$ws.Range("H132").Value = 3516.4783
$ws.Range("I132").Value = 3359.75
$ws.Range("J132").Value = 3874.7144
$ws.Range("K132").Value = 10079.25
$ws.Range("L132").Value = 11624.1432
$ws.Range("M132").Value = -7549.25
$ws.Range("N132").Value = -16684.1432

$ws.Range("H137").Value = 6210.625
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents() | Out-Null

$ws.Range("H138").Value = 4190.443
$ws.Range("J138").Value = 5589.2
$ws.Range("L138").Value = 16767.6
$ws.Range("N138").Value = -27047.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 702.1053000000001
$ws.Range("I5").Value = 943.0714
$ws.Range("J5").Value = 27.4
$ws.Range("K5").Value = 943.0714
$ws.Range("L5").Value = 27.4
$ws.Range("M5").Value = -831.0714
$ws.Range("N5").Value = -251.4

$ws.Range("H61").Value = 4212.385
$ws.Range("I61").Value = 4167
$ws.Range("K61").Value = 4167
$ws.Range("M61").Value = -3955

$ws.Range("H88").Value = 4984.5
$ws.Range("J88").Value = 4969
$ws.Range("L88").Value = 4969
$ws.Range("N88").Value = -5781

$ws.Range("H91").Value = 4984.5
$ws.Range("J91").Value = 4969
$ws.Range("L91").Value = 4969
$ws.Range("N91").Value = -7777

$ws.Range("H132").Value = 4239.05
$ws.Range("I132").Value = 4142.641
$ws.Range("K132").Value = 12427.923
$ws.Range("M132").Value = -9897.922999999999

$ws.Range("H136").Value = 4212.385
$ws.Range("I136").Value = 4167
$ws.Range("K136").Value = 12501
$ws.Range("M136").Value = -9951

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 702.1053000000001
$ws.Range("I4").Value = 943.0714
$ws.Range("J4").Value = 27.4
$ws.Range("K4").Value = 943.0714
$ws.Range("L4").Value = 27.4
$ws.Range("M4").Value = -828.0714
$ws.Range("N4").Value = -257.4

$ws.Range("H86").Value = 3165.5
$ws.Range("I86").Value = 2641.5
$ws.Range("J86").Value = 4999.5
$ws.Range("K86").Value = 2641.5
$ws.Range("L86").Value = 4999.5
$ws.Range("M86").Value = -1518.5
$ws.Range("N86").Value = -7245.5

$ws.Range("H89").Value = 3165.5
$ws.Range("I89").Value = 2641.5
$ws.Range("J89").Value = 4999.5
$ws.Range("K89").Value = 13207.5
$ws.Range("L89").Value = 24997.5
$ws.Range("M89").Value = -7591.5
$ws.Range("N89").Value = -36229.5

$ws.Range("H134").Value = 3114.077
$ws.Range("I134").Value = 2392.476
$ws.Range("K134").Value = 7177.428
$ws.Range("M134").Value = -4642.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5332.9
$ws.Range("I58").Value = 5544.625
$ws.Range("K58").Value = 5544.625
$ws.Range("M58").Value = -5341.625

$ws.Range("H122").Value = 2806.3333
$ws.Range("I122").Value = 1972
$ws.Range("K122").Value = 5916
$ws.Range("M122").Value = -3466

$ws.Range("H132").Value = 4530.2
$ws.Range("I132").Value = 3544.125
$ws.Range("K132").Value = 10632.375
$ws.Range("M132").Value = -8102.375

$ws.Range("H136").Value = 5332.9
$ws.Range("I136").Value = 5544.625
$ws.Range("K136").Value = 16633.875
$ws.Range("M136").Value = -14083.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 5461.3335
$ws.Range("I11").Value = 6611.25
$ws.Range("K11").Value = 19833.75
$ws.Range("M11").Value = -19693.75

$ws.Range("H69").Value = 3662.5
$ws.Range("J69").Value = 2985
$ws.Range("L69").Value = 8955
$ws.Range("N69").Value = -10577

$ws.Range("H72").Value = 3662.5
$ws.Range("J72").Value = 2985
$ws.Range("L72").Value = 26865
$ws.Range("N72").Value = -34977

$ws.Range("H121").Value = 167438.17
$ws.Range("I121").Value = 247.25
$ws.Range("J121").Value = 251033.62
$ws.Range("K121").Value = 741.75
$ws.Range("L121").Value = 753100.86
$ws.Range("M121").Value = 568.25
$ws.Range("N121").Value = -755720.86

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5591.0386
$ws.Range("I122").Value = 3844
$ws.Range("K122").Value = 11532
$ws.Range("M122").Value = -9082

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9461.467000000001
$ws.Range("I40").Value = 3778.5715
$ws.Range("K40").Value = 3778.5715
$ws.Range("M40").Value = -3642.5715

$ws.Range("H82").Value = 959.2143
$ws.Range("I82").Value = 959.2222
$ws.Range("J82").Value = 959.2
$ws.Range("K82").Value = 959.2222
$ws.Range("L82").Value = 959.2
$ws.Range("M82").Value = -598.2222
$ws.Range("N82").Value = -1681.2

$ws.Range("H85").Value = 959.2143
$ws.Range("I85").Value = 959.2222
$ws.Range("J85").Value = 959.2
$ws.Range("K85").Value = 959.2222
$ws.Range("L85").Value = 959.2
$ws.Range("M85").Value = 288.7778
$ws.Range("N85").Value = -3455.2

$ws.Range("H122").Value = 3931.4375
$ws.Range("I122").Value = 3726.8667
$ws.Range("K122").Value = 11180.6001
$ws.Range("M122").Value = -8730.6001

$ws.Range("H132").Value = 21576.648
$ws.Range("I132").Value = 28465.777
$ws.Range("K132").Value = 85397.33099999999
$ws.Range("M132").Value = -82867.33099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1846.421
$ws.Range("I81").Value = 2503.818
$ws.Range("J81").Value = 942.5
$ws.Range("K81").Value = 5007.636
$ws.Range("L81").Value = 1885
$ws.Range("M81").Value = -3946.636
$ws.Range("N81").Value = -4007

$ws.Range("H84").Value = 1846.421
$ws.Range("I84").Value = 2503.818
$ws.Range("J84").Value = 942.5
$ws.Range("K84").Value = 25038.18
$ws.Range("L84").Value = 9425
$ws.Range("M84").Value = -19734.18
$ws.Range("N84").Value = -20033

$ws.Range("H126").Value = 4781.24
$ws.Range("I126").Value = 4696.4
$ws.Range("J126").Value = 5120.6
$ws.Range("K126").Value = 14089.2
$ws.Range("L126").Value = 15361.8
$ws.Range("M126").Value = -11619.2
$ws.Range("N126").Value = -20301.8

$ws.Range("H132").Value = 1335.7084
$ws.Range("I132").Value = 1242.9474
$ws.Range("K132").Value = 3728.8422
$ws.Range("M132").Value = -1198.8422

$ws.Range("H136").Value = 1336.625
$ws.Range("I136").Value = 1410.8182
$ws.Range("K136").Value = 4232.4546
$ws.Range("M136").Value = -1682.4546
